$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lunsj")

# Fix typos in the "Beskrivelse" column (C2:C6):
#  - "Serveres" -> "Serverast"
#  - "hjemmelaget" -> "heimelaga"
#  - (row 4 only) "vårløk" -> "vårlauk"
$ws.Range("C2").Value = "Serverast med salat og heimelaga dressing. 1,2, 8"
$ws.Range("C3").Value = "Omelett med skinke, ost og vårløk. Serverast med salat og heimelaga dressing. 2, 3, 8"
$ws.Range("C4").Value = "Omelett med paprika, vårlauk og ost. Serverast med salat og heimelaga dressing. 2, 3, 8"
$ws.Range("C5").Value = "Omelett med ost og bacon. Serverast med salat og heimelaga dressing. 2, 3, 8"
$ws.Range("C6").Value = "Pai med ost, skinke, vårløk og fløte. Serverast med salat og heimelaga dressing. 1, 2, 3, 8"

# Update the selection to match the saved state of the workbook (active cell C6)
$ws.Range("C6").Select()

$wb.Save()
